# Weekly update to the "Choclo" (corn) price sheet:
# Two new daily records are inserted at the top of the data block
# (rows 266-267), pushing all subsequent records down by two rows.
#
# This is modeled in Excel as inserting two whole rows at row 266 and
# then filling those two new rows with the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 266 - everything that was on
# rows 266-353 moves down to rows 268-355.
$ws.Rows("266:267").Insert()

# Row 266: new "Choclero" record
$ws.Range("A266").Value = 4
$ws.Range("B266").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C266").Value = "Los Lagos"
$ws.Range("D266").Value = 44988
$ws.Range("E266").Value = 10
$ws.Range("F266").Value = 100112024
$ws.Range("G266").Value = "Choclo"
$ws.Range("H266").Value = "Choclero"
$ws.Range("I266").Value = "Primera"
$ws.Range("J266").Value = 5000
$ws.Range("K266").Value = 600
$ws.Range("L266").Value = 600
$ws.Range("M266").Value = 600
$ws.Range("N266").Value = "`$/unidad"
$ws.Range("O266").Value = "Región del Maule"
$ws.Range("P266").Value = 600
$ws.Range("Q266").Value = 1
$ws.Range("R266").Value = "Hortaliza"

# Row 267: new "Dulce o Americano" record
$ws.Range("A267").Value = 4
$ws.Range("B267").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C267").Value = "Los Lagos"
$ws.Range("D267").Value = 44988
$ws.Range("E267").Value = 10
$ws.Range("F267").Value = 100112024
$ws.Range("G267").Value = "Choclo"
$ws.Range("H267").Value = "Dulce o Americano"
$ws.Range("I267").Value = "Primera"
$ws.Range("J267").Value = 15000
$ws.Range("K267").Value = 250
$ws.Range("L267").Value = 250
$ws.Range("M267").Value = 250
$ws.Range("N267").Value = "`$/unidad"
$ws.Range("O267").Value = "Región de La Araucanía"
$ws.Range("P267").Value = 250
$ws.Range("Q267").Value = 1
$ws.Range("R267").Value = "Hortaliza"
